$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.445.66"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "1.616.76"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D12").Value = "1.844.19"
$ws.Range("E12").Value = "  +1.53%  "

$ws.Range("D13").Value = "1.619.14"
$ws.Range("E13").Value = "  +1.80%  "

$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "236.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.91%  "

$ws.Range("D18").Value = "26.441.90"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.66%  "

$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("E23").Value = "  +4.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.73%  "

$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("D32").Value = "1.516.29"
$ws.Range("E32").Value = "  +6.68%  "

$ws.Range("E33").Value = "  +1.55%  "

$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("E35").Value = "  +4.89%  "

$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("E38").Value = "  +0.24%  "

$ws.Range("E39").Value = "  +0.56%  "

$ws.Range("E40").Value = "  +2.22%  "

$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D43").Value = "1.756.18"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.913"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.77%  "

$ws.Range("E48").Value = "  -2.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("E50").Value = "  +0.14%  "

$ws.Range("E51").Value = "  +1.03%  "
